$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stratified_Data")

# Row 19: Source / Florida Health Department (repeated across B:Z)
$ws.Range("A19").Value = "Source"
$ws.Range("B19:Z19").Value = "Florida Health Department"

# Row 20: URL / testurl.org (repeated across B:Z)
$ws.Range("A20").Value = "URL"
$ws.Range("B20:Z20").Value = "testurl.org"

# Match the bold style used by column A elsewhere in the sheet
$ws.Range("A19").Style = $ws.Range("A18").Style
$ws.Range("A20").Style = $ws.Range("A18").Style

# Update selection to match the new active cell/selection recorded in the workbook
$ws.Range("C20:Z20").Select
$excel.ActiveWindow.RangeSelection.Item(1).Activate

# Remove the stale selection on the Comments sheet (reset to default view)
$commentsSheet = $wb.Worksheets.Item("Comments")
$commentsSheet.Range("A1").Select
